$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update results for previously pending rows (48, 49, 67, 68, 70)
$ws.Range("G48").Value = "Fallo"
$ws.Range("H48").Value = -1

$ws.Range("G49").Value = "Acierto"
$ws.Range("H49").Value = 1.3

$ws.Range("G67").Value = "Fallo"
$ws.Range("H67").Value = -1

$ws.Range("G68").Value = "Fallo"
$ws.Range("H68").Value = -1

$ws.Range("G70").Value = "Acierto"
$ws.Range("H70").Value = 1.25

# Append new tracked event as row 74
$ws.Range("A74").Value = 14310246

# Force column B to keep the date as literal text, not an auto-converted date serial
$ws.Range("B74").NumberFormat = "@"
$ws.Range("B74").Value = "2025-08-06"
$ws.Range("B74").Style = "Normal"

$ws.Range("C74").Value = "Mats Rosenkranz"
$ws.Range("D74").Value = "Lautaro Midon"
$ws.Range("E74").Value = "Gana Mats Rosenkranz"
$ws.Range("F74").Value = 3

# Result/profit are still pending for this newly tracked match, so leave them
# blank, but materialize the cells (matching the empty placeholder cells used
# for other still-pending rows such as G73/H73).
$ws.Range("G74").NumberFormat = "@"
$ws.Range("G74").Style = "Normal"
$ws.Range("H74").NumberFormat = "@"
$ws.Range("H74").Style = "Normal"
